$d = $word.ActiveDocument

# The two <id>...</id> runs (p160r_1 and p160r_2) were each split across three
# separate runs: "<id>", the bare identifier, and "</id>". Collapse each trio
# back into a single run by doing a literal Find/Replace over the full text
# span - Word merges multi-run matches into one run carrying the formatting
# of the first run in the match (Courier New / color 7f6000 / sz 18), which
# is exactly the formatting the surviving run should end up with.
# NOTE: "fig_p160r_1"/"fig_p160r_2" ids must stay untouched (and they will,
# since the search text below does not match them - the "fig_" prefix runs
# are a different paragraph and don't equal "<id>p160r_N</id>").

$d.Content.Find.Execute("<id>p160r_1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p160r_1</id>", 2) | Out-Null
$d.Content.Find.Execute("<id>p160r_2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p160r_2</id>", 2) | Out-Null
